# Add the new trade row (row 10) to the HZNP named-trade sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (date style on A, boolean style on G) from the row above
# so the new row's styles match the rest of the table exactly.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("G9").Copy()
$ws.Range("G10").PasteSpecial(-4122)

# Populate the new trade values.
$ws.Range("A10").Value = 42654.743842592594
$ws.Range("B10").Value = $true
$ws.Range("C10").Value = 10185.17
$ws.Range("D10").Value = 10012.950000000001
$ws.Range("E10").Value = 18.870000999999998
$ws.Range("F10").Value = 19.52
$ws.Range("G10").Value = $false
$ws.Range("H10").Value = 3.44
$ws.Range("I10").Value = $false

# Re-fit the column widths now that the table has an extra row of data
# (mirrors the bestFit column width recalculation Excel performs when the
# sheet's data changes).
$ws.Columns.Item(1).ColumnWidth = 14.5
$ws.Columns.Item(2).ColumnWidth = 7.5
$ws.Columns.Item(3).ColumnWidth = 8
$ws.Columns.Item(4).ColumnWidth = 10.5
$ws.Columns.Item(5).ColumnWidth = 9
$ws.Columns.Item(6).ColumnWidth = 6.1666666666666666
$ws.Columns.Item(7).ColumnWidth = 9.5
$ws.Columns.Item(8).ColumnWidth = 13.833333333333334
$ws.Columns.Item(9).ColumnWidth = 11
